$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.75159999999997
$ws.Range("A7").Value = -20.30759999999999
$ws.Range("D7").Value = -7.105199999999997
$ws.Range("D15").Value = -8.474199999999994
$ws.Range("A16").Value = -21.66159999999999
$ws.Range("D21").Value = -8.585899999999995
$ws.Range("D22").Value = -7.939400000000006
$ws.Range("D23").Value = -7.201799999999999
$ws.Range("A28").Value = -22.249
$ws.Range("A29").Value = -21.09239999999997
$ws.Range("A32").Value = -21.11209999999999
$ws.Range("D34").Value = -8.016999999999998
$ws.Range("A40").Value = -20.10289999999999
$ws.Range("D43").Value = -8.233300000000002
$ws.Range("D45").Value = -7.753699999999998
$ws.Range("D50").Value = -8.072699999999999
$ws.Range("D51").Value = -7.492499999999996
$ws.Range("A52").Value = -22.15649999999999
$ws.Range("A57").Value = -22.32
$ws.Range("A66").Value = -21.36239999999999
$ws.Range("D66").Value = -7.515000000000006
$ws.Range("D67").Value = -6.504700000000001
$ws.Range("D79").Value = -6.282100000000002
$ws.Range("D84").Value = -8.774400000000004
$ws.Range("D92").Value = -6.550600000000003
$ws.Range("D97").Value = -8.620500000000003
$ws.Range("A100").Value = -21.95319999999999
